$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit right after
#    "Yme van der Graaf" (just before the page break into the TOC section).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Locate the "Opdrachtbeschrijving" Heading 1 paragraph.
# ---------------------------------------------------------------------------
$headingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    $paraText = $p.Range.Text.TrimEnd()
    if ($styleName -eq "Heading 1" -and $paraText -eq "Opdrachtbeschrijving") {
        $headingIndex = $i
    }
}

# Right after the heading there are two empty paragraphs before the page
# break that leads into the "Planning" section. Keep the first one empty,
# and use the second one as the insertion point for the new body text plus
# the new "_GoBack" bookmark paragraph.
$targetPara = $d.Paragraphs.Item($headingIndex + 2)

$newBodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Het is de bedoeling dat er programma's worden geschreven die een goederenoverslag simuleren. Hiervoor worden XML bestanden ingelezen, die aangeven wanneer een schip, trein of vrachtwagen het terrein bereikt heeft, begint met laden en lossen, hoeveel containers er op het voertuig zitten en wanneer hij klaar is met laden en lossen. De verschillende containers die van het voertuig afkomen, worden vervolgens opgeslagen op het opslagterrein en anderen meteen over gedragen naar een ander voertuig, die de container weer verder vervoerd.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Op het terrein staan kranen die de containers van de voertuigen afhalen, dit zijn railkranen. Zij verplaatsen de containers één voor één naar de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>AGV's</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, die ze naar hun volgende bestemming rijden. Dit kan het opslagterrein zijn, waar ze met een vrij beweegbare kraan op het terrein worden gestapeld, of een volgend voertuig, waar ze met een railkraan op het voertuig worden gestapeld.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>De data moet uiteindelijk ook weergegeven worden op een mobiel apparaat.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="9000" w:name="_GoBack"/><w:bookmarkEnd w:id="9000"/></w:p>
'@

$targetPara.Range.InsertXML($newBodyXml)

Write-Output "ok"
